$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the data set (delete the
# lower-numbered row second so indices above it aren't disturbed).
$ws.Rows.Item(28).Delete()  # "SC 92" row
$ws.Rows.Item(26).Delete()  # "RM 232" row

# Individual cell edits (values that flipped between populated / missing).
$ws.Range("E5").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("D19").Value = -15.5
$ws.Range("E19").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("E25").Value = -7.1

# After the row deletions, the remaining rows have shifted up; re-apply
# the edits that land on the now-shifted rows.
$ws.Range("F26").Value = 17.38
$ws.Range("D27").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("D33").Value = -14.1
